$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet name: "Through 2021-10-24" -> "Through 2021-10-25"
$ws.Name = "Through 2021-10-25"

# Header label for the current (partial) month column
$ws.Range("B1").Value = "October 2021 (through October 25)"

# Numeric cell updates (existing cells whose values increased by 1)
$ws.Range("AF2").Value = 6
$ws.Range("B5").Value = 2
$ws.Range("L6").Value = 2
$ws.Range("BJ6").Value = 2
$ws.Range("AP10").Value = 4
$ws.Range("AF13").Value = 3
$ws.Range("B14").Value = 7
$ws.Range("L19").Value = 3
$ws.Range("L20").Value = 3
$ws.Range("AP32").Value = 3
$ws.Range("L36").Value = 7
$ws.Range("AP45").Value = 2
$ws.Range("L50").Value = 3
$ws.Range("AZ63").Value = 3
$ws.Range("B97").Value = 2

# Newly populated cells (were previously blank)
$ws.Range("B20").Value = 1
$ws.Range("V34").Value = 1
$ws.Range("V46").Value = 1
$ws.Range("B53").Value = 1
$ws.Range("B64").Value = 1
$ws.Range("AP70").Value = 1
$ws.Range("AF76").Value = 1
$ws.Range("BJ93").Value = 1
